# Auto-generated edit script for cryptos.xlsx update
# Updates Price (D) and Volume(1h) (E) columns, and for rows that were
# reordered (23/24, 32/33, 40/41) also updates Coin (B) and Link (C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.147.16'
$ws.Range("E2").Value = '  +2.48%  '

# Row 3
$ws.Range("D3").Value = '3.448.41'
$ws.Range("E3").Value = '  +1.62%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").Value = '''580.88'
$ws.Range("E5").Value = '  +3.88%  '

# Row 6
$ws.Range("D6").Value = '''187.19'
$ws.Range("E6").Value = '  +6.62%  '

# Row 7
$ws.Range("D7").Value = '''0.629'
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").Value = '3.440.60'
$ws.Range("E8").Value = '  +1.61%  '

# Row 9
$ws.Range("E9").Value = '  +0.03%  '

# Row 10
$ws.Range("E10").Value = '  -1.31%  '

# Row 11
$ws.Range("D11").Value = '''0.644'
$ws.Range("E11").Value = '  +1.01%  '

# Row 12
$ws.Range("D12").Value = '''57.61'
$ws.Range("E12").Value = '  +7.09%  '

# Row 13
$ws.Range("D13").Value = '''0.0000276'
$ws.Range("E13").Value = '  -1.57%  '

# Row 14
$ws.Range("D14").Value = '''9.46'
$ws.Range("E14").Value = '  +2.79%  '

# Row 15
$ws.Range("D15").Value = '3.986.53'
$ws.Range("E15").Value = '  +1.56%  '

# Row 16
$ws.Range("D16").Value = '''18.99'
$ws.Range("E16").Value = '  +3.56%  '

# Row 17
$ws.Range("D17").Value = '3.452.38'
$ws.Range("E17").Value = '  +2.13%  '

# Row 18
$ws.Range("D18").Value = '67.033.03'
$ws.Range("E18").Value = '  +2.72%  '

# Row 19
$ws.Range("E19").Value = '  -0.38%  '

# Row 20
$ws.Range("D20").Value = '''12.06'
$ws.Range("E20").Value = '  +1.65%  '

# Row 21
$ws.Range("E21").Value = '  +2.01%  '

# Row 22
$ws.Range("D22").Value = '''480.94'
$ws.Range("E22").Value = '  +4.19%  '

# Row 23
$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").Value = '''17.31'
$ws.Range("E23").Value = '  +22.52%  '

# Row 24
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '''5.33'
$ws.Range("E24").Value = '  +9.20%  '

# Row 25
$ws.Range("D25").Value = '''4.35'
$ws.Range("E25").Value = '  +5.77%  '

# Row 26
$ws.Range("D26").Value = '''89.28'
$ws.Range("E26").Value = '  +1.73%  '

# Row 27
$ws.Range("E27").Value = '  +1.82%  '

# Row 28
$ws.Range("D28").Value = '''10.97'
$ws.Range("E28").Value = '  +2.45%  '

# Row 29
$ws.Range("D29").Value = '''9.02'
$ws.Range("E29").Value = '  +3.07%  '

# Row 30
$ws.Range("D30").Value = '''31.17'
$ws.Range("E30").Value = '  -0.14%  '

# Row 31
$ws.Range("D31").Value = '''7.35'
$ws.Range("E31").Value = '  +12.03%  '

# Row 32
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '''602.54'
$ws.Range("E32").Value = '  +3.71%  '

# Row 33
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").Value = '''64.81'
$ws.Range("E33").Value = '  +2.05%  '

# Row 34
$ws.Range("D34").Value = '''11.79'
$ws.Range("E34").Value = '  +2.50%  '

# Row 35
$ws.Range("D35").Value = '''0.112'
$ws.Range("E35").Value = '  +3.51%  '

# Row 36
$ws.Range("E36").Value = '  -0.07%  '

# Row 37
$ws.Range("D37").Value = '''0.147'
$ws.Range("E37").Value = '  +2.74%  '

# Row 38
$ws.Range("D38").Value = '''36.94'
$ws.Range("E38").Value = '  +3.10%  '

# Row 39
$ws.Range("D39").Value = '''0.388'
$ws.Range("E39").Value = '  +3.64%  '

# Row 40
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0765'
$ws.Range("E40").Value = '  +2.67%  '

# Row 41
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '''3.48'
$ws.Range("E41").Value = '  -4.22%  '

# Row 42
$ws.Range("D42").Value = '3.198.48'
$ws.Range("E42").Value = '  +3.29%  '

# Row 43
$ws.Range("D43").Value = '''2.89'
$ws.Range("E43").Value = '  +3.95%  '

# Row 44
$ws.Range("D44").Value = '''0.0429'
$ws.Range("E44").Value = '  +2.23%  '

# Row 45
$ws.Range("D45").Value = '''2.59'
$ws.Range("E45").Value = '  +5.69%  '

# Row 46
$ws.Range("D46").Value = '''3.24'
$ws.Range("E46").Value = '  +1.82%  '

# Row 47
$ws.Range("D47").Value = '''0.135'
$ws.Range("E47").Value = '  +0.93%  '

# Row 48
$ws.Range("D48").Value = '''2.70'
$ws.Range("E48").Value = '  +18.09%  '

# Row 49
$ws.Range("D49").Value = '''0.999'
$ws.Range("E49").Value = '  +0.16%  '

# Row 50
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").Value = '''8.64'
$ws.Range("E50").Value = '  +3.28%  '

# Row 51
$ws.Range("D51").Value = '''3.21'
$ws.Range("E51").Value = '  +3.85%  '

